# chore(runtime): publish files + archive (2025-12-07 15:08:00)
#
# 1) Refresh the "scraped_at" timestamps (column K, rows 2-35) on the
#    "snapshot" sheet with the latest scrape run's timestamps.
# 2) Remove the now-resolved entry (row 2) from the "new_injured" sheet,
#    shrinking its used range back down to just the header row.

$wb = $excel.ActiveWorkbook

# --- 1. Update scraped_at timestamps on the "snapshot" sheet -----------
$snapshot = $wb.Worksheets.Item("snapshot")

$scrapedAt = @(
    "2025-12-07T07:02:11.936156+00:00",
    "2025-12-07T07:02:11.936200+00:00",
    "2025-12-07T07:02:11.936224+00:00",
    "2025-12-07T07:02:14.733530+00:00",
    "2025-12-07T07:02:14.733564+00:00",
    "2025-12-07T07:02:16.968234+00:00",
    "2025-12-07T07:02:19.768546+00:00",
    "2025-12-07T07:02:22.103783+00:00",
    "2025-12-07T07:02:24.974303+00:00",
    "2025-12-07T07:02:31.263968+00:00",
    "2025-12-07T07:02:31.264001+00:00",
    "2025-12-07T07:02:33.767145+00:00",
    "2025-12-07T07:02:36.279350+00:00",
    "2025-12-07T07:02:38.630391+00:00",
    "2025-12-07T07:02:40.934179+00:00",
    "2025-12-07T07:02:40.934212+00:00",
    "2025-12-07T07:02:48.494941+00:00",
    "2025-12-07T07:02:48.494981+00:00",
    "2025-12-07T07:02:48.495002+00:00",
    "2025-12-07T07:02:50.708967+00:00",
    "2025-12-07T07:02:50.708997+00:00",
    "2025-12-07T07:02:53.481160+00:00",
    "2025-12-07T07:02:53.481195+00:00",
    "2025-12-07T07:02:53.481216+00:00",
    "2025-12-07T07:02:57.129901+00:00",
    "2025-12-07T07:03:06.085269+00:00",
    "2025-12-07T07:03:06.085301+00:00",
    "2025-12-07T07:03:06.085320+00:00",
    "2025-12-07T07:03:08.400943+00:00",
    "2025-12-07T07:03:08.400970+00:00",
    "2025-12-07T07:03:11.212443+00:00",
    "2025-12-07T07:03:11.212472+00:00",
    "2025-12-07T07:03:13.976446+00:00",
    "2025-12-07T07:03:13.976477+00:00"
)

$startRow = 2
for ($i = 0; $i -lt $scrapedAt.Length; $i++) {
    $row = $startRow + $i
    $snapshot.Cells.Item($row, 11).Value = $scrapedAt[$i]
}

# --- 2. Drop the resolved row from the "new_injured" sheet -------------
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()
